$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update footer timestamp text (row 1)
$ws.Range("A1").Value = "Datos actualizados a 23 de Marzo de 2020 a las 23:16"

# Update country names and statistics that changed due to reordering / data refresh
# Row 4
$ws.Range("C4").Value = 0
$ws.Range("G4").Value = 0

# Row 6
$ws.Range("B6").Value = 43022
$ws.Range("C6").Value = 9456
$ws.Range("E6").Value = 42197
$ws.Range("G6").Value = 117
$ws.Range("H6").Value = 530

# Row 7
$ws.Range("B7").Value = 35136
$ws.Range("C7").Value = 6368
$ws.Range("E7").Value = 29470
$ws.Range("G7").Value = 539
$ws.Range("H7").Value = 2311

# Row 18
$ws.Range("A18").Value = "Canada"
$ws.Range("B18").Value = 2091
$ws.Range("C18").Value = 621
$ws.Range("D18").Value = 320
$ws.Range("E18").Value = 1748
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 3

# Row 19
$ws.Range("A19").Value = "Portugal"
$ws.Range("B19").Value = 2060
$ws.Range("C19").Value = 460
$ws.Range("D19").Value = 14
$ws.Range("E19").Value = 2023
$ws.Range("F19").Value = 47
$ws.Range("G19").Value = 9

# Row 73
$ws.Range("A73").Value = "Lituania"
$ws.Range("B73").Value = 179
$ws.Range("D73").Value = 1
$ws.Range("E73").Value = 177
$ws.Range("F73").Value = 1
$ws.Range("G73").Value = 0
$ws.Range("H73").Value = 1

# Row 74
$ws.Range("A74").Value = "Hungria"
$ws.Range("B74").Value = 167
$ws.Range("C74").Value = 36
$ws.Range("D74").Value = 16
$ws.Range("E74").Value = 144
$ws.Range("F74").Value = 6
$ws.Range("G74").Value = 1
$ws.Range("H74").Value = 7

# Row 112
$ws.Range("A112").Value = "Ruanda"
$ws.Range("C112").Value = 17
$ws.Range("E112").Value = 36
$ws.Range("H112").Value = 0

# Row 113
$ws.Range("A113").Value = "Consejo Danes para los Refugiados"
$ws.Range("C113").Value = 6
$ws.Range("E113").Value = 35
$ws.Range("F113").Value = 0
$ws.Range("H113").Value = 1

# Row 114
$ws.Range("A114").Value = "Mauricio"
$ws.Range("C114").Value = 8
$ws.Range("D114").Value = 0
$ws.Range("E114").Value = 34
$ws.Range("F114").Value = 1
$ws.Range("G114").Value = 0
$ws.Range("H114").Value = 2

# Row 115
$ws.Range("A115").Value = "Nigeria"
$ws.Range("B115").Value = 36
$ws.Range("D115").Value = 2
$ws.Range("E115").Value = 33
$ws.Range("H115").Value = 1

# Row 116
$ws.Range("A116").Value = "Banglades"
$ws.Range("B116").Value = 33
$ws.Range("C116").Value = 6
$ws.Range("D116").Value = 5
$ws.Range("E116").Value = 25
$ws.Range("H116").Value = 3

# Row 117
$ws.Range("A117").Value = "Puerto Rico"
$ws.Range("B117").Value = 31
$ws.Range("C117").Value = 8
$ws.Range("D117").Value = 1
$ws.Range("G117").Value = 1
$ws.Range("H117").Value = 2

# Row 118
$ws.Range("A118").Value = "Guam"
$ws.Range("B118").Value = 29
$ws.Range("C118").Value = 2
$ws.Range("E118").Value = 28
$ws.Range("H118").Value = 1

# Row 120
$ws.Range("A120").Value = "Honduras"
$ws.Range("C120").Value = 1
$ws.Range("E120").Value = 27
$ws.Range("G120").Value = 0
$ws.Range("H120").Value = 0

# Row 121
$ws.Range("A121").Value = "Montenegro"
$ws.Range("C121").Value = 6
$ws.Range("E121").Value = 26
$ws.Range("H121").Value = 1

# Row 122
$ws.Range("A122").Value = "Ghana"
$ws.Range("B122").Value = 27
$ws.Range("C122").Value = 4
$ws.Range("D122").Value = 0
$ws.Range("E122").Value = 25
$ws.Range("G122").Value = 1
$ws.Range("H122").Value = 2

# Row 123
$ws.Range("A123").Value = "Costa de Marfil"
$ws.Range("C123").Value = 11
$ws.Range("D123").Value = 2
$ws.Range("E123").Value = 23

# Row 124
$ws.Range("A124").Value = "Macao"
$ws.Range("B124").Value = 25
$ws.Range("C124").Value = 3
$ws.Range("D124").Value = 10
$ws.Range("E124").Value = 15

# Row 125
$ws.Range("A125").Value = "Mayotte"
$ws.Range("B125").Value = 24
$ws.Range("C125").Value = 13
$ws.Range("D125").Value = 0
$ws.Range("E125").Value = 24

# Row 126
$ws.Range("A126").Value = "Monaco"
$ws.Range("B126").Value = 23
$ws.Range("D126").Value = 1
$ws.Range("E126").Value = 22
$ws.Range("F126").Value = 0
$ws.Range("H126").Value = 0

# Row 127
$ws.Range("A127").Value = "Paraguay"
$ws.Range("B127").Value = 22
$ws.Range("C127").Value = 0
$ws.Range("E127").Value = 21
$ws.Range("F127").Value = 1

# Row 129
$ws.Range("A129").Value = "Guyana"
$ws.Range("C129").Value = 1
$ws.Range("D129").Value = 0
$ws.Range("E129").Value = 19
$ws.Range("H129").Value = 1

# Row 130
$ws.Range("A130").Value = "Guayana Francesa"
$ws.Range("B130").Value = 20
$ws.Range("C130").Value = 2
$ws.Range("D130").Value = 6
$ws.Range("E130").Value = 14

# Row 132
$ws.Range("A132").Value = "Togo"
$ws.Range("C132").Value = 2

# Row 133
$ws.Range("A133").Value = "Polinesia Francesa"
$ws.Range("C133").Value = 0

# Row 136
$ws.Range("A136").Value = "Kirguistan"
$ws.Range("C136").Value = 2

# Row 137
$ws.Range("A137").Value = "Kenia"
$ws.Range("C137").Value = 1

# Row 153
$ws.Range("A153").Value = "Benin"
$ws.Range("C153").Value = 3

# Row 154
$ws.Range("A154").Value = "Surinam"
$ws.Range("C154").Value = 0
$ws.Range("E154").Value = 5
$ws.Range("H154").Value = 0

# Row 155
$ws.Range("A155").Value = "Islas Caimanes"
$ws.Range("C155").Value = 2

# Row 156
$ws.Range("A156").Value = "Gabon"
$ws.Range("B156").Value = 5
$ws.Range("H156").Value = 1

# Row 157
$ws.Range("A157").Value = "Namibia"
$ws.Range("C157").Value = 1

# Row 158
$ws.Range("A158").Value = "Bahamas"

# Row 160
$ws.Range("A160").Value = "Suazilandia"

# Row 161
$ws.Range("A161").Value = "Guinea"
$ws.Range("C161").Value = 2

# Row 162
$ws.Range("A162").Value = "Groenlandia"
$ws.Range("C162").Value = 0
$ws.Range("E162").Value = 4
$ws.Range("H162").Value = 0

# Row 163
$ws.Range("A163").Value = "Curazao"
$ws.Range("B163").Value = 4
$ws.Range("C163").Value = 1
$ws.Range("H163").Value = 1

# Row 164
$ws.Range("A164").Value = "Antigua y Barbuda"
$ws.Range("C164").Value = 2

# Row 165
$ws.Range("A165").Value = "Fiyi"
$ws.Range("C165").Value = 1

# Row 166
$ws.Range("A166").Value = "Santa Lucia"
$ws.Range("C166").Value = 1

# Row 167
$ws.Range("A167").Value = "Republica de Africa Central"

# Row 168
$ws.Range("A168").Value = "Zambia"

# Row 169
$ws.Range("A169").Value = "El Salvador"

# Row 170
$ws.Range("A170").Value = "Cabo Verde"
$ws.Range("C170").Value = 0

# Row 171
$ws.Range("A171").Value = "Republica de Yibuti"
$ws.Range("C171").Value = 2

# Row 172
$ws.Range("A172").Value = "San Bartolome"
$ws.Range("C172").Value = 0

# Row 173
$ws.Range("A173").Value = "Angola"
$ws.Range("C173").Value = 1

# Row 174
$ws.Range("A174").Value = "Liberia"
$ws.Range("E174").Value = 3
$ws.Range("G174").Value = 0
$ws.Range("H174").Value = 0

# Row 175
$ws.Range("A175").Value = "Zimbabue"
$ws.Range("B175").Value = 3
$ws.Range("G175").Value = 1
$ws.Range("H175").Value = 1

# Row 176
$ws.Range("A176").Value = "Nicaragua"

# Row 177
$ws.Range("A177").Value = "Butan"
$ws.Range("C177").Value = 0

# Row 178
$ws.Range("A178").Value = "San Martin (Parte Holandesa)"
$ws.Range("C178").Value = 1

# Row 179
$ws.Range("A179").Value = "Birmania"
$ws.Range("C179").Value = 2
